$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a single date serial (45180) repeated for
# every data row (rows 2-189). The workbook was refreshed one day later,
# so every one of those cells advances by exactly one day (45180 -> 45181).
$ws.Range("C2:C189").Value = 45181
